{"js": "// Prefix the five \"Heading 3\" section titles with their new numeric\n// labels, mirroring the PR-45 preview deploy edit:\n//   Description -> 0.1 Description\n//   Usage       -> 0.2 Usage\n//   Arguments   -> 0.3 Arguments\n//   Value       -> 0.4 Value\n//   Examples    -> 0.5 Examples\n\nconst numberedHeadings = {\n  \"Description\": \"0.1 Description\",\n  \"Usage\": \"0.2 Usage\",\n  \"Arguments\": \"0.3 Arguments\",\n  \"Value\": \"0.4 Value\",\n  \"Examples\": \"0.5 Examples\"\n};\n\nconst body = context.document.body;\nconst paragraphs = body.paragraphs;\nparagraphs.load(\"items/text,items/style\");\nawait context.sync();\n\nfor (let i = 0; i < paragraphs.items.length; i++) {\n  const paragraph = paragraphs.items[i];\n  const newText = numberedHeadings[paragraph.text];\n  if (paragraph.style === \"Heading 3\" && newText) {\n    paragraph.insertText(newText, Word.InsertLocation.replace);\n  }\n}\n\nawait context.sync();\n", "ps1": "# Prefix the five \"Heading 3\" section titles with their new numeric\n# labels, mirroring the PR-45 preview deploy edit:\n#   Description -> 0.1 Description\n#   Usage       -> 0.2 Usage\n#   Arguments   -> 0.3 Arguments\n#   Value       -> 0.4 Value\n#   Examples    -> 0.5 Examples\n\n$d = $word.ActiveDocument\n\n$numberedHeadings = @{\n  \"Description\" = \"0.1 Description\"\n  \"Usage\"       = \"0.2 Usage\"\n  \"Arguments\"   = \"0.3 Arguments\"\n  \"Value\"       = \"0.4 Value\"\n  \"Examples\"    = \"0.5 Examples\"\n}\n\nforeach ($paragraph in $d.Paragraphs) {\n  $range = $paragraph.Range\n  $plainText = $range.Text.TrimEnd(\"`r\")\n  if ($range.Style.NameLocal -eq \"Heading 3\" -and $numberedHeadings.ContainsKey($plainText)) {\n    # Exclude the trailing paragraph mark so only the heading text is replaced.\n    $range.MoveEnd(1, -1) | Out-Null\n    $range.Text = $numberedHeadings[$plainText]\n  }\n}\n"}
